$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Target change (per the diff): the paragraph that carries the lone
# "_GoBack" bookmark gains a new "Customer" run placed right before the
# bookmark tags, and the two paragraphs that follow it ("User" and
# "List all users") are removed outright.
# ------------------------------------------------------------------

# Locate the bookmark paragraph defensively by scanning for the empty
# paragraph that is immediately followed by "User" / "List all users",
# rather than trusting a hard-coded index.
$bookmarkIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count - 2; $i++) {
    $cur   = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r")
    $next  = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd("`r")
    $next2 = $d.Paragraphs.Item($i + 2).Range.Text.TrimEnd("`r")
    if ($cur -eq "" -and $next -eq "User" -and $next2 -eq "List all users") {
        $bookmarkIndex = $i
        break
    }
}
if ($bookmarkIndex -eq -1) {
    throw "could not locate the bookmark / User / List all users paragraphs"
}

# Insert the "Customer" run in front of whatever is already in that
# paragraph (the bookmarkStart/bookmarkEnd pair), then stamp it with the
# same en-US run-language formatting used throughout the document.
$bookmarkPara = $d.Paragraphs.Item($bookmarkIndex)
$bookmarkPara.Range.InsertBefore("Customer")

$bookmarkPara = $d.Paragraphs.Item($bookmarkIndex)
$bookmarkPara.Range.LanguageID = "en-US"

# Remove the "User" paragraph entirely.
$userPara = $d.Paragraphs.Item($bookmarkIndex + 1)
$userText = $userPara.Range.Text.TrimEnd("`r")
if ($userText -ne "User") {
    throw "expected 'User' paragraph, found: $userText"
}
$userPara.Range.Delete()

# Remove the "List all users" paragraph entirely (it is now shifted into
# the same index since the "User" paragraph above it was just removed).
$listAllUsersPara = $d.Paragraphs.Item($bookmarkIndex + 1)
$listAllUsersText = $listAllUsersPara.Range.Text.TrimEnd("`r")
if ($listAllUsersText -ne "List all users") {
    throw "expected 'List all users' paragraph, found: $listAllUsersText"
}
$listAllUsersPara.Range.Delete()
